# Applies the "Add templates for staging, logging, sender, inboundProcessing"
# edit described by the diff:
#   1. Split the release-date text "2017/08/10" into four runs whose
#      concatenated text reads "2017/09/01".
#   2. Replace the trailing (_GoBack bookmark-only) paragraph with a new
#      "Object Inventory" Heading1 section, an "Integration Directory"
#      Heading2 section, a two-column "Configuration Scenario" table, and
#      a trailing empty paragraph.

$d = $word.ActiveDocument

# --- 1. Release date: "2017/08/10" -> runs "2017/0" + "9" + "/" + "01" ---
$dateParaText = "2017/08/10"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "$dateParaText*") {
        $dateRange = $candidate.Range
        break
    }
}

$dateXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>2017/0</w:t></w:r><w:r><w:t>9</w:t></w:r><w:r><w:t>/</w:t></w:r><w:r><w:t>01</w:t></w:r></w:p>'
[void]$dateRange.InsertXML($dateXml)

# --- 2. Replace the bookmark-only paragraph with the new sections ---
$bookmarkPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    $xml = $candidate.Range.WordOpenXML
    if ($xml -like "*_GoBack*") {
        $bookmarkPara = $candidate
        break
    }
}

$newSectionsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Obje</w:t></w:r><w:r><w:t>ct Inventory</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Integration Directory</w:t></w:r></w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblW w:w="5000" w:type="pct"/>
    <w:tblLook w:val="00A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="0"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1639"/>
    <w:gridCol w:w="7705"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="5000" w:type="pct"/>
        <w:gridSpan w:val="2"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
        </w:tcBorders>
        <w:shd w:val="clear" w:color="auto" w:fill="8496B0" w:themeFill="text2" w:themeFillTint="99"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:b/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
          </w:rPr>
          <w:t>Configuration Scenario</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="877" w:type="pct"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
        </w:tcBorders>
        <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:b/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
          </w:rPr>
          <w:t>Name</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4123" w:type="pct"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:keepNext/>
          <w:spacing w:line="276" w:lineRule="auto"/>
        </w:pPr>
        <w:r>
          <w:t>$</w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>ICO_Value</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:trPr>
      <w:trHeight w:val="345"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="877" w:type="pct"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
        </w:tcBorders>
        <w:shd w:val="clear" w:color="auto" w:fill="BFBFBF" w:themeFill="background1" w:themeFillShade="BF"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:b/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
          </w:rPr>
          <w:t>Description</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4123" w:type="pct"/>
        <w:tcBorders>
          <w:top w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:left w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:bottom w:val="single" w:sz="6" w:space="0" w:color="auto"/>
          <w:right w:val="single" w:sz="6" w:space="0" w:color="auto"/>
        </w:tcBorders>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:keepNext/>
          <w:spacing w:line="276" w:lineRule="auto"/>
          <w:rPr>
            <w:highlight w:val="yellow"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:t>$</w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>DESCRIPTION_Value</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

[void]$bookmarkPara.Range.InsertXML($newSectionsXml)
